# Insert a new data row before the current row 17 (shifts existing rows
# 17..149 down to 18..150, extending the used range to A1:R150) and
# populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("17").Insert()

$ws.Range("A17").Value = 5
$ws.Range("B17").Value = "Macroferia Regional de Talca"
$ws.Range("C17").Value = "Maule"
$ws.Range("D17").Value = 44670
$ws.Range("E17").Value = 7
$ws.Range("F17").Value = 100112030
$ws.Range("G17").Value = "Poroto granado"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 200
$ws.Range("K17").Value = 20000
$ws.Range("L17").Value = 20000
$ws.Range("M17").Value = 20000
$ws.Range("N17").Value = "`$/saco 25 kilos"
$ws.Range("O17").Value = "Región del Maule"
$ws.Range("P17").Value = 800
$ws.Range("Q17").Value = 25
$ws.Range("R17").Value = "Hortaliza"
